# EditShareSkill() and DeleteShareSkill() added.
# Update the "ShareSkill" test-data sheet: drop the "Available days" column,
# re-point the Skill Trade value, and add two new sample rows (an edited
# share-skill row and a hidden/deleted one) with wrapped, top-aligned text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShareSkill")

# --- 1. Remove the "Available days" column (old column J) -----------------
# Everything to its right (Start time, End time, Skill Trade, Skill-Exchange,
# Work Samples, Active) shifts one column left automatically.
$ws.Columns.Item(10).Delete()

# --- 2. Fix up the existing sample row (row 2) -----------------------------
# After the column shift, L2 still holds the old M2 value ("Credit"); the
# real "Skill Trade" value for this row is "Skill-exchanges".
$ws.Range("L2").Value = "Skill-exchanges"

# Row 2 becomes a wrapped / top-aligned, 60pt-tall row (matches the new
# "Quality Assurance" row added below).
$ws.Range("A2:O2").WrapText = $true
$ws.Range("A2:O2").VerticalAlignment = -4160
$ws.Rows.Item(2).RowHeight = 60

# --- 3. Add the new "Quality Assurance" share-skill row (row 3) -----------
$ws.Range("A3").Value = "Quality Assurance"
$ws.Range("B3").Value = "Code test automation scripts with a specialization in C Sharp, Selenium and other tools."
$ws.Range("C3").Value = "Business"
$ws.Range("F3").Value = "One-off service"
$ws.Range("L3").Value = "Credit"
$ws.Range("M3").Value = 10
$ws.Range("O3").Value = "Hidden"

$ws.Range("A3:O3").WrapText = $true
$ws.Range("A3:O3").VerticalAlignment = -4160
$ws.Rows.Item(3).RowHeight = 60

# B3's description keeps the same un-highlighted look as the other body
# cells above it (no row-level fill override needed here).

# --- 4. Blank placeholder time cells (rows 3-5, Start/End time columns) ---
$ws.Range("J3:K5").NumberFormat = "h:mm AM/PM"
$ws.Range("J3:K5").VerticalAlignment = -4160

# --- 5. Refresh the view so the selection matches the authored state ------
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("K3").Select()
